$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.184.45"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "2.594.01"
$ws.Range("E3").Value = "  +3.10%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'315.64"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").Value = "'97.70"
$ws.Range("E6").Value = "  +3.64%  "
$ws.Range("D7").Value = "'0.577"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.537"
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("D10").Value = "'35.77"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "2.992.94"
$ws.Range("E13").Value = "  +3.10%  "
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.565.22"
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'15.30"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("D17").Value = "'0.848"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "43.285.65"
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("D19").Value = "'6.85"
$ws.Range("E19").Value = "  +2.82%  "
$ws.Range("D20").Value = "'12.75"
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("D22").Value = "'69.66"
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("D23").Value = "'252.84"
$ws.Range("E23").Value = "  +0.77%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("E25").Value = "  +3.23%  "
$ws.Range("D26").Value = "'27.38"
$ws.Range("E26").Value = "  +2.43%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("D29").Value = "'41.02"
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("D30").Value = "'10.31"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("D31").Value = "'5.87"
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("D32").Value = "'156.48"
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("E33").Value = "  +5.81%  "
$ws.Range("D34").Value = "'2.16"
$ws.Range("E34").Value = "  +2.34%  "
$ws.Range("D35").Value = "'0.0807"
$ws.Range("E35").Value = "  +3.66%  "
$ws.Range("E36").Value = "  +3.08%  "
$ws.Range("D37").Value = "'18.80"
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("E38").Value = "  +2.17%  "
$ws.Range("E39").Value = "  +10.15%  "
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("D41").Value = "'23.29"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("D42").Value = "'3.98"
$ws.Range("E42").Value = "  +5.66%  "
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.018.19"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'3.24"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("E47").Value = "  +1.36%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "2.843.82"
$ws.Range("E48").Value = "  +3.17%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").Value = "'83.33"
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.196"
$ws.Range("E50").Value = "  +4.58%  "
$ws.Range("D51").Value = "'104.56"
$ws.Range("E51").Value = "  +2.46%  "
